# This sheet is a weekly "Fruta / hortaliza" price log for
# "Feria Lagunitas de Puerto Montt - Coliflor". A new weekly record needs to
# be inserted as a new row at position 593, which pushes all the following
# rows (593-636) down by one (to 594-637), growing the used range from
# A1:R636 to A1:R637.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 593; everything from the old row 593 onward
# shifts down by one row (old row 636 becomes row 637).
$ws.Rows("593:593").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A593").Value = 4
$ws.Range("B593").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C593").Value = "Los Lagos"
$ws.Range("D593").Value = 45265
$ws.Range("E593").Value = 10
$ws.Range("F593").Value = 100112008
$ws.Range("G593").Value = "Coliflor"
$ws.Range("H593").Value = "Sin especificar"
$ws.Range("I593").Value = "Primera"
$ws.Range("J593").Value = 1500
$ws.Range("K593").Value = 1600
$ws.Range("L593").Value = 1700
$ws.Range("M593").Value = 1650
$ws.Range("N593").Value = "$/unidad"
$ws.Range("O593").Value = "Región Metropolitana"
$ws.Range("P593").Value = 1650
$ws.Range("Q593").Value = 1
$ws.Range("R593").Value = "Hortaliza"
